$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.000.41'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '2.606.64'
$ws.Range('E3').Value = '  -1.63%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.54'
$ws.Range('E5').Value = '  +3.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.79'
$ws.Range('E6').Value = '  -1.97%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.596'
$ws.Range('E8').Value = '  +3.88%  '
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('E10').Value = '  -1.58%  '
$ws.Range('E11').Value = '  +5.34%  '
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('D13').Value = '3.064.73'
$ws.Range('E13').Value = '  -1.84%  '
$ws.Range('D14').Value = '58.962.44'
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.94'
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('D16').Value = '2.622.31'
$ws.Range('E16').Value = '  -1.84%  '
$ws.Range('E17').Value = '  -1.88%  '
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '338.79'
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('E20').Value = '  -2.15%  '
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.61'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  +2.41%  '
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.161'
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('E27').Value = '  -1.65%  '
$ws.Range('D28').Value = '0.0₃0757'
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  +1.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.99'
$ws.Range('E31').Value = '  +1.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '153.86'
$ws.Range('E32').Value = '  +2.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.00'
$ws.Range('E33').Value = '  +0.52%  '
$ws.Range('E34').Value = '  -1.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.897'
$ws.Range('E35').Value = '  +6.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.887'
$ws.Range('E36').Value = '  +5.28%  '
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.98'
$ws.Range('E38').Value = '  -1.00%  '
$ws.Range('E39').Value = '  +0.65%  '
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '284.22'
$ws.Range('E41').Value = '  -0.64%  '
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.602'
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0955'
$ws.Range('E44').Value = '  +0.66%  '
$ws.Range('E45').Value = '  -0.39%  '
$ws.Range('E46').Value = '  -0.97%  '
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.68'
$ws.Range('E48').Value = '  +2.30%  '
$ws.Range('D49').Value = '1.949.74'
$ws.Range('E49').Value = '  -0.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '117.43'
$ws.Range('E50').Value = '  +4.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.08'
$ws.Range('E51').Value = '  -2.06%  '
